# Applies the "New simulation files for schemes report" edit:
#  - replaces the list of HKL orientation/scheme names used down column B
#  - the per-HKL column header labels (row 2, C:M) get reshuffled to a new order
#  - 10 new simulation rows (20-29) are appended, each filled with 1's like the
#    existing data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels for row 2, columns C..W -----------------------------
# (N2:W2 keep the same text as before; C2:M2 are reordered)
$headerCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
$headerVals = @("[4, 2, 2]","[5, 1, 1]","[2, 2, 2]","[1, 1, 1]","[3, 1, 1]","[3, 3, 1]","[2, 2, 0]","[2, 0, 0]","[3, 3, 3]","[4, 0, 0]","[4, 2, 0]","1Pair-A","1Pair-B","2Pairs-A","2Pairs-B","3Pairs-A","3Pairs-B","3Pairs-C","4Pairs","5A4F","MaxUnique")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "2").Value = $headerVals[$i]
}

# --- New scheme names for column B, rows 3-29 ------------------------------
# Rows 3-19 replace the old scheme list; rows 20-29 are brand new rows.
$schemeNames = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

$dataCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $i + 3
    $ws.Range("B" + $row).Value = $schemeNames[$i]
}

# --- Append brand-new rows 20-29 (A/B already set above) ------------------
for ($row = 20; $row -le 29; $row++) {
    $ws.Range("A" + $row).Value = $row - 2

    foreach ($col in $dataCols) {
        $ws.Range($col + $row).Value = 1
    }
}

# Column A on the data rows carries the bold/bordered "s=1" style used by the
# header row/column; copy that formatting down onto the newly added rows.
$ws.Range("A19").Copy()
$ws.Range("A20:A29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
